$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '27.676.15'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  -0.15%  '
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.632.63'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  +0.15%  '
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '212.06'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  +0.11%  '
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '23.17'
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  -0.33%  '
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.264'
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  +1.62%  '
$ws.Range("E10").Value = '  +0.13%  '
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0863'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  -3.00%  '
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '1.865.83'
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  -0.35%  '
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '1.642.56'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  +0.05%  '
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '4.04'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  -0.35%  '
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '0.552'
$dCell.Style = "Normal"
$ws.Range("E15").Value = '  -1.83%  '
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '65.10'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  +0.47%  '
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '27.652.53'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  -0.14%  '
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '229.90'
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  -0.30%  '
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₃0719'
$dCell.Style = "Normal"
$ws.Range("E19").Value = '  -0.51%  '
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '7.57'
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("E21").Value = '  +0.11%  '
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '10.71'
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  +4.11%  '
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '4.34'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  +0.60%  '
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '2.13'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  +2.86%  '
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '148.74'
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  -1.50%  '
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '6.87'
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  -1.11%  '
$ws.Range("E27").Value = '  -1.02%  '
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '15.58'
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("E31").Value = '  -0.86%  '
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '3.29'
$dCell.Style = "Normal"
$ws.Range("E32").Value = '  -0.85%  '
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '1.471.65'
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("E34").Value = '  -1.16%  '
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '1.54'
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  -1.72%  '
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '2.33'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("E37").Value = '  +6.03%  '
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.877'
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("E39").Value = '  -1.73%  '
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("E41").Value = '  +0.78%  '
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '67.80'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  -1.83%  '
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '2.48'
$dCell.Style = "Normal"
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("E45").Value = '  -5.03%  '
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '1.774.44'
$dCell.Style = "Normal"
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("E47").Value = '  +0.43%  '
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '87.46'
$dCell.Style = "Normal"
$ws.Range("E48").Value = '  +0.55%  '
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₆0106'
$dCell.Style = "Normal"
$ws.Range("E49").Value = '  -1.00%  '
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0991'
$dCell.Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '7.74'
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  -1.23%  '
